$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = 2.65
$ws.Range("K2").Value = 2.2
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 8
$ws.Range("S2").Value = 1.35
$ws.Range("T2").Value = 2.95
$ws.Range("U2").Value = 1.6
$ws.Range("AC2").Value = 8
$ws.Range("AD2").Value = 6.9
$ws.Range("AQ2").Value = 40
$ws.Range("AT2").Value = 2.95
$ws.Range("AV2").Value = 55
$ws.Range("AX2").Value = 17

$wb.Save()
